$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain text so values like "59.987.09" or "0.998"
# are not reinterpreted as numbers/dates by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.987.09'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.408.81'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.71'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.79'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.584'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.44%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.62'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.350'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.63'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.837.49'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.879.41'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.408.99'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '326.59'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.73'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.68'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.70%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.59'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.41'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +5.14%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0765'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.53'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.12'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.10'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +8.31%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.40'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.00%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.46%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.19'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '321.93'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '145.59'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.63%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.90'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.87%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0513'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.575'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.06'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.67'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.86%  '
